$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46; this shifts existing rows 46:145 down to 47:146
# and carries formatting (e.g. the date style on column D) down with them.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new record.
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = 44544
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = 100112021
$ws.Range("G46").Value = "Ají"
$ws.Range("H46").Value = "Americana (o)"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 150
$ws.Range("K46").Value = 16000
$ws.Range("L46").Value = 16000
$ws.Range("M46").Value = 16000
$ws.Range("N46").Value = "$/caja 15 kilos"
$ws.Range("O46").Value = "Región del Maule"
$ws.Range("P46").Value = 1067
$ws.Range("Q46").Value = 15
$ws.Range("R46").Value = "Hortaliza"
